# ----------------------------------------------------------------------------
# Scheduled runner update: refresh Universalis market-price-derived columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) on every crafting-class sheet (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR) using the latest pulled prices.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3453.1765
$ws.Range("I98").Value = 1940.8518
$ws.Range("J98").Value = 9286.429
$ws.Range("K98").Value = 1940.8518
$ws.Range("L98").Value = 9286.429
$ws.Range("M98").Value = -442.8517999999999
$ws.Range("N98").Value = -12282.429
$ws.Range("H107").Value = 1863.25
$ws.Range("I107").Value = 1842.8572
$ws.Range("J107").Value = 2006
$ws.Range("K107").Value = 1842.8572
$ws.Range("L107").Value = 2006
$ws.Range("M107").Value = 77.14280000000008
$ws.Range("N107").Value = -5846
$ws.Range("H113").Value = 22000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 22000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 22000
$ws.Range("N113").Value = -28508
$ws.Range("H122").Value = 3453.1765
$ws.Range("I122").Value = 1940.8518
$ws.Range("J122").Value = 9286.429
$ws.Range("K122").Value = 5822.555399999999
$ws.Range("L122").Value = 27859.287
$ws.Range("M122").Value = -3372.555399999999
$ws.Range("N122").Value = -32759.287
$ws.Range("H138").Value = 3914.8867
$ws.Range("I138").Value = 1403.7222
$ws.Range("J138").Value = 4487.051
$ws.Range("K138").Value = 4211.1666
$ws.Range("L138").Value = 13461.153
$ws.Range("M138").Value = 928.8334000000004
$ws.Range("N138").Value = -23741.153
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1635.1666
$ws.Range("I2").Value = 1622.2
$ws.Range("J2").Value = 1700
$ws.Range("K2").Value = 1622.2
$ws.Range("L2").Value = 1700
$ws.Range("M2").Value = -1509.2
$ws.Range("N2").Value = -1926
$ws.Range("H116").Value = 1635.1666
$ws.Range("I116").Value = 1622.2
$ws.Range("J116").Value = 1700
$ws.Range("K116").Value = 1622.2
$ws.Range("L116").Value = 1700
$ws.Range("M116").Value = 671.8
$ws.Range("N116").Value = -6288
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H132").Value = 3251.7666
$ws.Range("I132").Value = 1487.0555
$ws.Range("J132").Value = 5898.8335
$ws.Range("K132").Value = 4461.166499999999
$ws.Range("L132").Value = 17696.5005
$ws.Range("M132").Value = -1931.166499999999
$ws.Range("N132").Value = -22756.5005
$ws.Range("H140").Value = 50429
$ws.Range("J140").Value = 50429
$ws.Range("L140").Value = 50429
$ws.Range("N140").Value = -60789
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1635.1666
$ws.Range("I3").Value = 1622.2
$ws.Range("J3").Value = 1700
$ws.Range("K3").Value = 1622.2
$ws.Range("L3").Value = 1700
$ws.Range("M3").Value = -1508.2
$ws.Range("N3").Value = -1928
$ws.Range("H134").Value = 2856.6614
$ws.Range("I134").Value = 1882.1904
$ws.Range("J134").Value = 4903.05
$ws.Range("K134").Value = 5646.5712
$ws.Range("L134").Value = 14709.15
$ws.Range("M134").Value = -3111.5712
$ws.Range("N134").Value = -19779.15

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 207.66667
$ws.Range("I7").Value = 131.28572
$ws.Range("J7").Value = 475
$ws.Range("K7").Value = 131.28572
$ws.Range("L7").Value = 475
$ws.Range("M7").Value = -18.28572
$ws.Range("N7").Value = -701
$ws.Range("H31").Value = 4223.2144
$ws.Range("I31").Value = 1178.0667
$ws.Range("J31").Value = 7736.846
$ws.Range("K31").Value = 1178.0667
$ws.Range("L31").Value = 7736.846
$ws.Range("M31").Value = -883.0667000000001
$ws.Range("N31").Value = -8326.846
$ws.Range("H34").Value = 4223.2144
$ws.Range("I34").Value = 1178.0667
$ws.Range("J34").Value = 7736.846
$ws.Range("K34").Value = 1178.0667
$ws.Range("L34").Value = 7736.846
$ws.Range("M34").Value = -976.0667000000001
$ws.Range("N34").Value = -8140.846
$ws.Range("H58").Value = 2348.6462
$ws.Range("I58").Value = 1732.6271
$ws.Range("J58").Value = 8406.166999999999
$ws.Range("K58").Value = 1732.6271
$ws.Range("L58").Value = 8406.166999999999
$ws.Range("M58").Value = -1529.6271
$ws.Range("N58").Value = -8812.166999999999
$ws.Range("H99").Value = 3669.1365
$ws.Range("I99").Value = 1981.4
$ws.Range("J99").Value = 7285.7144
$ws.Range("K99").Value = 1981.4
$ws.Range("L99").Value = 7285.7144
$ws.Range("M99").Value = -483.4000000000001
$ws.Range("N99").Value = -10281.7144
$ws.Range("H126").Value = 3669.1365
$ws.Range("I126").Value = 1981.4
$ws.Range("J126").Value = 7285.7144
$ws.Range("K126").Value = 5944.200000000001
$ws.Range("L126").Value = 21857.1432
$ws.Range("M126").Value = -3474.200000000001
$ws.Range("N126").Value = -26797.1432
$ws.Range("H132").Value = 3090.0977
$ws.Range("I132").Value = 2441.0588
$ws.Range("J132").Value = 6242.5713
$ws.Range("K132").Value = 7323.176399999999
$ws.Range("L132").Value = 18727.7139
$ws.Range("M132").Value = -4793.176399999999
$ws.Range("N132").Value = -23787.7139
$ws.Range("H136").Value = 2348.6462
$ws.Range("I136").Value = 1732.6271
$ws.Range("J136").Value = 8406.166999999999
$ws.Range("K136").Value = 5197.8813
$ws.Range("L136").Value = 25218.501
$ws.Range("M136").Value = -2647.8813
$ws.Range("N136").Value = -30318.501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 13520170
$ws.Range("J131").Value = 1110
$ws.Range("L131").Value = 3330
$ws.Range("N131").Value = -13410

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3466.138
$ws.Range("I102").Value = 2948.16
$ws.Range("K102").Value = 2948.16
$ws.Range("M102").Value = -1326.16
$ws.Range("H132").Value = 3094.9
$ws.Range("I132").Value = 1630.8334
$ws.Range("J132").Value = 4070.9443
$ws.Range("K132").Value = 4892.5002
$ws.Range("L132").Value = 12212.8329
$ws.Range("M132").Value = -2362.5002
$ws.Range("N132").Value = -17272.8329
$ws.Range("H139").Value = 66500
$ws.Range("J139").Value = 66500
$ws.Range("L139").Value = 66500
$ws.Range("N139").Value = -76780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1415.2354
$ws.Range("I61").Value = 1394.9231
$ws.Range("J61").Value = 1481.25
$ws.Range("K61").Value = 1394.9231
$ws.Range("L61").Value = 1481.25
$ws.Range("M61").Value = -1192.9231
$ws.Range("N61").Value = -1885.25
$ws.Range("H113").Value = 1415.2354
$ws.Range("I113").Value = 1394.9231
$ws.Range("J113").Value = 1481.25
$ws.Range("K113").Value = 1394.9231
$ws.Range("L113").Value = 1481.25
$ws.Range("M113").Value = 775.0769
$ws.Range("N113").Value = -5821.25
$ws.Range("H122").Value = 3089.6897
$ws.Range("I122").Value = 2725.963
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 8177.889000000001
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -5727.889000000001
$ws.Range("N122").Value = -28900
$ws.Range("H132").Value = 8359.048000000001
$ws.Range("I132").Value = 3059
$ws.Range("J132").Value = 11620.615
$ws.Range("K132").Value = 9177
$ws.Range("L132").Value = 34861.845
$ws.Range("M132").Value = -6647
$ws.Range("N132").Value = -39921.845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 27819
$ws.Range("I58").Value = 22442.25
$ws.Range("J58").Value = 34988
$ws.Range("K58").Value = 22442.25
$ws.Range("L58").Value = 34988
$ws.Range("M58").Value = -22134.25
$ws.Range("N58").Value = -35604
$ws.Range("I62").Value = 250002500
$ws.Range("J62").Value = 308750.75
$ws.Range("K62").Value = 250002500
$ws.Range("L62").Value = 308750.75
$ws.Range("M62").Value = -250001876
$ws.Range("N62").Value = -309998.75
$ws.Range("I65").Value = 250002500
$ws.Range("J65").Value = 308750.75
$ws.Range("K65").Value = 1250012500
$ws.Range("L65").Value = 1543753.75
$ws.Range("M65").Value = -1250009380
$ws.Range("N65").Value = -1549993.75
$ws.Range("H132").Value = 4274973.5
$ws.Range("I132").Value = 850.1129
$ws.Range("K132").Value = 2550.3387
$ws.Range("M132").Value = -20.33869999999979
$ws.Range("H138").Value = 52671.145
$ws.Range("J138").Value = 52671.145
$ws.Range("L138").Value = 52671.145
$ws.Range("N138").Value = -62951.145
